$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.335.15"
$ws.Range("E2").Value = "  -2.18%  "

# Row 3
$ws.Range("D3").Value = "1.646.28"
$ws.Range("E3").Value = "  -3.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").Value = "'310.40"
$ws.Range("E5").Value = "  -0.61%  "

# Row 7
$ws.Range("D7").Value = "'0.3641"
$ws.Range("E7").Value = "  -3.16%  "

# Row 8
$ws.Range("D8").Value = "'46.68"
$ws.Range("E8").Value = "  -5.99%  "

# Row 9
$ws.Range("D9").Value = "'0.3223"
$ws.Range("E9").Value = "  -6.73%  "

# Row 10
$ws.Range("D10").Value = "'1.111"
$ws.Range("E10").Value = "  -8.34%  "

# Row 11
$ws.Range("D11").Value = "'0.06978"
$ws.Range("E11").Value = "  -7.79%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.19%  "

# Row 13
$ws.Range("D13").Value = "'5.897"
$ws.Range("E13").Value = "  -7.11%  "

# Row 14
$ws.Range("D14").Value = "'19.19"
$ws.Range("E14").Value = "  -9.94%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.644.72"
$ws.Range("E15").Value = "  -3.70%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'6.536"
$ws.Range("E16").Value = "  -7.77%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.06570"
$ws.Range("E17").Value = "  -2.33%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.00001026"
$ws.Range("E18").Value = "  -9.51%  "

# Row 20
$ws.Range("D20").Value = "'77.57"
$ws.Range("E20").Value = "  -9.08%  "

# Row 21
$ws.Range("D21").Value = "'5.881"
$ws.Range("E21").Value = "  -8.29%  "

# Row 22
$ws.Range("D22").Value = "'15.44"
$ws.Range("E22").Value = "  -11.28%  "

# Row 23
$ws.Range("D23").Value = "'12.40"
$ws.Range("E23").Value = "  -7.02%  "

# Row 24
$ws.Range("D24").Value = "24.343.67"
$ws.Range("E24").Value = "  -2.17%  "

# Row 25
$ws.Range("D25").Value = "'2.474"
$ws.Range("E25").Value = "  +0.72%  "

# Row 26
$ws.Range("D26").Value = "'2.270"
$ws.Range("E26").Value = "  -19.23%  "

# Row 27
$ws.Range("D27").Value = "'145.34"
$ws.Range("E27").Value = "  -4.21%  "

# Row 28
$ws.Range("E28").Value = "  -10.40%  "

# Row 29
$ws.Range("D29").Value = "1.827.07"
$ws.Range("E29").Value = "  -3.81%  "

# Row 30
$ws.Range("D30").Value = "'123.21"
$ws.Range("E30").Value = "  -7.48%  "

# Row 31
$ws.Range("D31").Value = "'1.157"
$ws.Range("E31").Value = "  -7.23%  "

# Row 32
$ws.Range("D32").Value = "'4.060"
$ws.Range("E32").Value = "  -4.28%  "

# Row 33
$ws.Range("D33").Value = "'5.580"
$ws.Range("E33").Value = "  -19.84%  "

# Row 34
$ws.Range("D34").Value = "'0.08399"
$ws.Range("E34").Value = "  -5.12%  "

# Row 35
$ws.Range("D35").Value = "'1.650"
$ws.Range("E35").Value = "  -7.99%  "

# Row 36
$ws.Range("D36").Value = "'11.92"
$ws.Range("E36").Value = "  -14.58%  "

# Row 37
$ws.Range("D37").Value = "'5.132"
$ws.Range("E37").Value = "  -8.99%  "

# Row 38
$ws.Range("D38").Value = "'1.242"
$ws.Range("E38").Value = "  -3.36%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05931"
$ws.Range("E39").Value = "  -11.34%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02196"
$ws.Range("E40").Value = "  -9.10%  "

# Row 41
$ws.Range("D41").Value = "'0.2034"
$ws.Range("E41").Value = "  -9.32%  "

# Row 42
$ws.Range("D42").Value = "'8.018"
$ws.Range("E42").Value = "  -14.54%  "

# Row 43
$ws.Range("E43").Value = "  +0.32%  "

# Row 44
$ws.Range("D44").Value = "'0.5820"
$ws.Range("E44").Value = "  -10.06%  "

# Row 45
$ws.Range("D45").Value = "'3.749"
$ws.Range("E45").Value = "  -2.17%  "

# Row 46
$ws.Range("D46").Value = "'12.49"
$ws.Range("E46").Value = "  -10.60%  "

# Row 47
$ws.Range("D47").Value = "'0.5533"
$ws.Range("E47").Value = "  -10.55%  "

# Row 48
$ws.Range("D48").Value = "'121.50"
$ws.Range("E48").Value = "  -6.94%  "

# Row 49
$ws.Range("D49").Value = "'1.923"
$ws.Range("E49").Value = "  -10.19%  "

# Row 50
$ws.Range("D50").Value = "'0.06868"
$ws.Range("E50").Value = "  -6.24%  "

# Row 51
$ws.Range("E51").Value = "  -4.37%  "
